$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix D685:D687 - drop the orphan "Arial 10 family2 (no charset)" font/style
# (cellXf 15 / font 7) in favour of the plain default font, as in the diff.
$ws.Range("D685:D687").Font.Name = "Arial"
$ws.Range("D685:D687").Font.Size = 10
$ws.Range("D685:D687").NumberFormat = "General"

# Row 695
$ws.Cells.Item(695,1).Value = 'Protéalpes'
$ws.Cells.Item(695,2).Value = 'B'
$ws.Cells.Item(695,3).Value = 'Boisson sucre'
$ws.Cells.Item(695,4).Value = 'Boisson Endurance Fruits Rouges ou Menthe'
$ws.Cells.Item(695,5).Value = 1
$ws.Cells.Item(695,6).Value = 91.5
$ws.Cells.Item(695,7).Formula = '=(E695/100)*F695'
$ws.Cells.Item(695,8).Value = 0
$ws.Cells.Item(695,9).Formula = '=(E695/100)*H695'
$ws.Cells.Item(695,10).Value = 0
$ws.Cells.Item(695,11).Formula = '=G695/E695'
$ws.Cells.Item(695,11).NumberFormat = "0.0000"
$ws.Cells.Item(695,12).Formula = '=0.576/100'
$ws.Cells.Item(695,13).Formula = '=(27.95/800)/G695'
$ws.Cells.Item(695,14).Value = 1
$ws.Cells.Item(695,15).Value = 0
$ws.Cells.Item(695,16).Value = 0
$ws.Cells.Item(695,17).Value = 0
$ws.Cells.Item(695,21).Value = 'NF EN 17444'
$ws.Cells.Item(695,22).Value = 1

# Row 696
$ws.Cells.Item(696,1).Value = 'Ultimum '
$ws.Cells.Item(696,2).Value = 'C'
$ws.Cells.Item(696,3).Value = 'Compote sucre'
$ws.Cells.Item(696,4).Value = 'Ultimum Sport Banane'
$ws.Cells.Item(696,5).Value = 70
$ws.Cells.Item(696,6).Value = 59
$ws.Cells.Item(696,7).Formula = '=(E696/100)*F696'
$ws.Cells.Item(696,8).Value = 0.6
$ws.Cells.Item(696,9).Formula = '=(E696/100)*H696'
$ws.Cells.Item(696,10).Value = 0
$ws.Cells.Item(696,11).Formula = '=G696/E696'
$ws.Cells.Item(696,11).NumberFormat = "0.0000"
$ws.Cells.Item(696,12).Value = 0
$ws.Cells.Item(696,13).Formula = '=2.5/G696'
$ws.Cells.Item(696,14).Value = 0
$ws.Cells.Item(696,15).Value = 0
$ws.Cells.Item(696,16).Value = 0
$ws.Cells.Item(696,17).Value = 0
$ws.Cells.Item(696,22).Value = 0

# Row 697
$ws.Cells.Item(697,1).Value = 'Ultimum '
$ws.Cells.Item(697,2).Value = 'C'
$ws.Cells.Item(697,3).Value = 'Compote sucre'
$ws.Cells.Item(697,4).Value = 'Ultimum Sport Figue'
$ws.Cells.Item(697,5).Value = 70
$ws.Cells.Item(697,6).Formula = '=52/0.7'
$ws.Cells.Item(697,7).Formula = '=(E697/100)*F697'
$ws.Cells.Item(697,8).Formula = '=1.1/0.7'
$ws.Cells.Item(697,9).Formula = '=(E697/100)*H697'
$ws.Cells.Item(697,10).Value = 0
$ws.Cells.Item(697,11).Formula = '=G697/E697'
$ws.Cells.Item(697,11).NumberFormat = "0.0000"
$ws.Cells.Item(697,12).Value = 0
$ws.Cells.Item(697,13).Formula = '=2.5/G697'
$ws.Cells.Item(697,14).Value = 0
$ws.Cells.Item(697,15).Value = 0
$ws.Cells.Item(697,16).Value = 0
$ws.Cells.Item(697,17).Value = 0
$ws.Cells.Item(697,22).Value = 0

# Row 698
$ws.Cells.Item(698,1).Value = 'Ultimum '
$ws.Cells.Item(698,2).Value = 'C'
$ws.Cells.Item(698,3).Value = 'Compote sucre'
$ws.Cells.Item(698,4).Value = 'Ultimum Sport Pruneau'
$ws.Cells.Item(698,5).Value = 70
$ws.Cells.Item(698,6).Formula = '=44/0.7'
$ws.Cells.Item(698,7).Formula = '=(E698/100)*F698'
$ws.Cells.Item(698,8).Formula = '=0.7/0.7'
$ws.Cells.Item(698,9).Formula = '=(E698/100)*H698'
$ws.Cells.Item(698,10).Value = 0
$ws.Cells.Item(698,11).Formula = '=G698/E698'
$ws.Cells.Item(698,11).NumberFormat = "0.0000"
$ws.Cells.Item(698,12).Value = 0
$ws.Cells.Item(698,13).Formula = '=2.5/G698'
$ws.Cells.Item(698,14).Value = 0
$ws.Cells.Item(698,15).Value = 0
$ws.Cells.Item(698,16).Value = 0
$ws.Cells.Item(698,17).Value = 0
$ws.Cells.Item(698,22).Value = 0

# Row 699
$ws.Cells.Item(699,1).Value = 'Ultimum '
$ws.Cells.Item(699,2).Value = 'C'
$ws.Cells.Item(699,3).Value = 'Compote sucre'
$ws.Cells.Item(699,4).Value = 'Ultimum Sport Abricot'
$ws.Cells.Item(699,5).Value = 70
$ws.Cells.Item(699,6).Formula = '=47/0.7'
$ws.Cells.Item(699,7).Formula = '=(E699/100)*F699'
$ws.Cells.Item(699,8).Formula = '=1.2/0.7'
$ws.Cells.Item(699,9).Formula = '=(E699/100)*H699'
$ws.Cells.Item(699,10).Value = 0
$ws.Cells.Item(699,11).Formula = '=G699/E699'
$ws.Cells.Item(699,11).NumberFormat = "0.0000"
$ws.Cells.Item(699,12).Value = 0
$ws.Cells.Item(699,13).Formula = '=2.5/G699'
$ws.Cells.Item(699,14).Value = 0
$ws.Cells.Item(699,15).Value = 0
$ws.Cells.Item(699,16).Value = 0
$ws.Cells.Item(699,17).Value = 0
$ws.Cells.Item(699,22).Value = 0

# Row 700
$ws.Cells.Item(700,1).Value = 'Ultimum '
$ws.Cells.Item(700,2).Value = 'C'
$ws.Cells.Item(700,3).Value = 'Compote sucre'
$ws.Cells.Item(700,4).Value = 'Ultimum Sport Date'
$ws.Cells.Item(700,5).Value = 70
$ws.Cells.Item(700,6).Formula = '=39/0.7'
$ws.Cells.Item(700,7).Formula = '=(E700/100)*F700'
$ws.Cells.Item(700,8).Formula = '=1/0.7'
$ws.Cells.Item(700,9).Formula = '=(E700/100)*H700'
$ws.Cells.Item(700,10).Value = 0
$ws.Cells.Item(700,11).Formula = '=G700/E700'
$ws.Cells.Item(700,11).NumberFormat = "0.0000"
$ws.Cells.Item(700,12).Value = 0
$ws.Cells.Item(700,13).Formula = '=2.5/G700'
$ws.Cells.Item(700,14).Value = 0
$ws.Cells.Item(700,15).Value = 0
$ws.Cells.Item(700,16).Value = 0
$ws.Cells.Item(700,17).Value = 0
$ws.Cells.Item(700,22).Value = 0

# Row 701
$ws.Cells.Item(701,1).Value = 'Ultimum '
$ws.Cells.Item(701,2).Value = 'C'
$ws.Cells.Item(701,3).Value = 'Compote sucre'
$ws.Cells.Item(701,4).Value = 'Ultimum Sport Ananas'
$ws.Cells.Item(701,5).Value = 70
$ws.Cells.Item(701,6).Formula = '=52/0.7'
$ws.Cells.Item(701,7).Formula = '=(E701/100)*F701'
$ws.Cells.Item(701,8).Formula = '=0.4/0.7'
$ws.Cells.Item(701,9).Formula = '=(E701/100)*H701'
$ws.Cells.Item(701,10).Value = 0
$ws.Cells.Item(701,11).Formula = '=G701/E701'
$ws.Cells.Item(701,11).NumberFormat = "0.0000"
$ws.Cells.Item(701,12).Formula = '=0.04*0.4'
$ws.Cells.Item(701,13).Formula = '=2.5/G701'
$ws.Cells.Item(701,14).Value = 0
$ws.Cells.Item(701,15).Value = 0
$ws.Cells.Item(701,16).Value = 0
$ws.Cells.Item(701,17).Value = 0
$ws.Cells.Item(701,22).Value = 0

# Row 702
$ws.Cells.Item(702,1).Value = 'Ultimum '
$ws.Cells.Item(702,2).Value = 'C'
$ws.Cells.Item(702,3).Value = 'Compote sucre'
$ws.Cells.Item(702,4).Value = 'Ultimum Sport Mangue'
$ws.Cells.Item(702,5).Value = 70
$ws.Cells.Item(702,6).Formula = '=54.7/0.7'
$ws.Cells.Item(702,7).Formula = '=(E702/100)*F702'
$ws.Cells.Item(702,8).Formula = '=0.8/0.7'
$ws.Cells.Item(702,9).Formula = '=(E702/100)*H702'
$ws.Cells.Item(702,10).Value = 0
$ws.Cells.Item(702,11).Formula = '=G702/E702'
$ws.Cells.Item(702,11).NumberFormat = "0.0000"
$ws.Cells.Item(702,12).Formula = '=0.35*0.4'
$ws.Cells.Item(702,13).Formula = '=2.5/G702'
$ws.Cells.Item(702,14).Value = 0
$ws.Cells.Item(702,15).Value = 0
$ws.Cells.Item(702,16).Value = 0
$ws.Cells.Item(702,17).Value = 0
$ws.Cells.Item(702,22).Value = 0

# Row 703
$ws.Cells.Item(703,1).Value = 'Ultimum '
$ws.Cells.Item(703,2).Value = 'C'
$ws.Cells.Item(703,3).Value = 'Compote sucre'
$ws.Cells.Item(703,4).Value = 'Ultimum Sport Kiwi ou Pêche  '
$ws.Cells.Item(703,5).Value = 70
$ws.Cells.Item(703,6).Formula = '=40/0.7'
$ws.Cells.Item(703,7).Formula = '=(E703/100)*F703'
$ws.Cells.Item(703,8).Formula = '=0.4/0.7'
$ws.Cells.Item(703,9).Formula = '=(E703/100)*H703'
$ws.Cells.Item(703,10).Value = 0
$ws.Cells.Item(703,11).Formula = '=G703/E703'
$ws.Cells.Item(703,11).NumberFormat = "0.0000"
$ws.Cells.Item(703,12).Value = 0.016
$ws.Cells.Item(703,13).Formula = '=2.5/G703'
$ws.Cells.Item(703,14).Value = 0
$ws.Cells.Item(703,15).Value = 0
$ws.Cells.Item(703,16).Value = 0
$ws.Cells.Item(703,17).Value = 0
$ws.Cells.Item(703,22).Value = 0

# Row 704
$ws.Cells.Item(704,1).Value = 'Ultimum '
$ws.Cells.Item(704,2).Value = 'C'
$ws.Cells.Item(704,3).Value = 'Compote sucre'
$ws.Cells.Item(704,4).Value = 'Ultimum Sport Mix Pruneau Cranberry ou Banane Cassis ou Date Goji'
$ws.Cells.Item(704,5).Value = 70
$ws.Cells.Item(704,6).Formula = '=40/0.7'
$ws.Cells.Item(704,7).Formula = '=(E704/100)*F704'
$ws.Cells.Item(704,8).Formula = '=0.5/0.7'
$ws.Cells.Item(704,9).Formula = '=(E704/100)*H704'
$ws.Cells.Item(704,10).Value = 0
$ws.Cells.Item(704,11).Formula = '=G704/E704'
$ws.Cells.Item(704,11).NumberFormat = "0.0000"
$ws.Cells.Item(704,12).Value = 0.01
$ws.Cells.Item(704,13).Formula = '=2.5/G704'
$ws.Cells.Item(704,14).Value = 0
$ws.Cells.Item(704,15).Value = 0
$ws.Cells.Item(704,16).Value = 0
$ws.Cells.Item(704,17).Value = 0
$ws.Cells.Item(704,22).Value = 0

# Row 705
$ws.Cells.Item(705,1).Value = 'Ultimum '
$ws.Cells.Item(705,2).Value = 'C'
$ws.Cells.Item(705,3).Value = 'Compote sucre'
$ws.Cells.Item(705,4).Value = 'Ultimum Oxygen Datte ou Pruneaux ou Abricot'
$ws.Cells.Item(705,5).Value = 70
$ws.Cells.Item(705,6).Formula = '=20.3/0.7'
$ws.Cells.Item(705,7).Formula = '=(E705/100)*F705'
$ws.Cells.Item(705,8).Formula = '=0.56/0.7'
$ws.Cells.Item(705,9).Formula = '=(E705/100)*H705'
$ws.Cells.Item(705,10).Value = 0
$ws.Cells.Item(705,11).Formula = '=G705/E705'
$ws.Cells.Item(705,11).NumberFormat = "0.0000"
$ws.Cells.Item(705,12).Value = 0
$ws.Cells.Item(705,13).Formula = '=2.5/G705'
$ws.Cells.Item(705,14).Value = 0
$ws.Cells.Item(705,15).Value = 0
$ws.Cells.Item(705,16).Value = 0
$ws.Cells.Item(705,17).Value = 0
$ws.Cells.Item(705,22).Value = 0

# Row 706
$ws.Cells.Item(706,1).Value = 'Ultimum '
$ws.Cells.Item(706,2).Value = 'C'
$ws.Cells.Item(706,3).Value = 'Compote sucre'
$ws.Cells.Item(706,4).Value = 'Ultimum Oxygen Figue'
$ws.Cells.Item(706,5).Value = 70
$ws.Cells.Item(706,6).Formula = '=16.59/0.7'
$ws.Cells.Item(706,7).Formula = '=(E706/100)*F706'
$ws.Cells.Item(706,8).Formula = '=1.05/0.7'
$ws.Cells.Item(706,9).Formula = '=(E706/100)*H706'
$ws.Cells.Item(706,10).Value = 0
$ws.Cells.Item(706,11).Formula = '=G706/E706'
$ws.Cells.Item(706,11).NumberFormat = "0.0000"
$ws.Cells.Item(706,12).Value = 0
$ws.Cells.Item(706,13).Formula = '=2.5/G706'
$ws.Cells.Item(706,14).Value = 0
$ws.Cells.Item(706,15).Value = 0
$ws.Cells.Item(706,16).Value = 0
$ws.Cells.Item(706,17).Value = 0
$ws.Cells.Item(706,22).Value = 0

# Row 707
$ws.Cells.Item(707,1).Value = 'Ultimum '
$ws.Cells.Item(707,2).Value = 'C'
$ws.Cells.Item(707,3).Value = 'Compote sucre'
$ws.Cells.Item(707,4).Value = 'Ultimum Oxygen Mangue'
$ws.Cells.Item(707,5).Value = 70
$ws.Cells.Item(707,6).Formula = '=25.3/0.7'
$ws.Cells.Item(707,7).Formula = '=(E707/100)*F707'
$ws.Cells.Item(707,8).Formula = '=0.8/0.7'
$ws.Cells.Item(707,9).Formula = '=(E707/100)*H707'
$ws.Cells.Item(707,10).Value = 0
$ws.Cells.Item(707,11).Formula = '=G707/E707'
$ws.Cells.Item(707,11).NumberFormat = "0.0000"
$ws.Cells.Item(707,12).Value = 0
$ws.Cells.Item(707,13).Formula = '=2.5/G707'
$ws.Cells.Item(707,14).Value = 0
$ws.Cells.Item(707,15).Value = 0
$ws.Cells.Item(707,16).Value = 0
$ws.Cells.Item(707,17).Value = 0
$ws.Cells.Item(707,22).Value = 0

# Row 708
$ws.Cells.Item(708,1).Value = 'Ultimum '
$ws.Cells.Item(708,2).Value = 'CS'
$ws.Cells.Item(708,3).Value = 'Compote sel'
$ws.Cells.Item(708,4).Value = 'Ultimum Mix Salé Petits Pois Kiwi'
$ws.Cells.Item(708,5).Value = 70
$ws.Cells.Item(708,6).Formula = '=6.72/0.7'
$ws.Cells.Item(708,7).Formula = '=(E708/100)*F708'
$ws.Cells.Item(708,8).Formula = '=2.52/0.7'
$ws.Cells.Item(708,9).Formula = '=(E708/100)*H708'
$ws.Cells.Item(708,10).Value = 0
$ws.Cells.Item(708,11).Formula = '=G708/E708'
$ws.Cells.Item(708,11).NumberFormat = "0.0000"
$ws.Cells.Item(708,12).Formula = '=0.42*0.4'
$ws.Cells.Item(708,13).Formula = '=2.5/G708'
$ws.Cells.Item(708,14).Value = 0
$ws.Cells.Item(708,15).Value = 0
$ws.Cells.Item(708,16).Value = 0
$ws.Cells.Item(708,17).Value = 0
$ws.Cells.Item(708,22).Value = 0

# Row 709
$ws.Cells.Item(709,1).Value = 'Ultimum '
$ws.Cells.Item(709,2).Value = 'CS'
$ws.Cells.Item(709,3).Value = 'Compote sel'
$ws.Cells.Item(709,4).Value = 'Ultimum Mix Salé Tomate Carotte'
$ws.Cells.Item(709,5).Value = 70
$ws.Cells.Item(709,6).Formula = '=3.01/0.7'
$ws.Cells.Item(709,7).Formula = '=(E709/100)*F709'
$ws.Cells.Item(709,8).Formula = '=0.63/0.7'
$ws.Cells.Item(709,9).Formula = '=(E709/100)*H709'
$ws.Cells.Item(709,10).Value = 0
$ws.Cells.Item(709,11).Formula = '=G709/E709'
$ws.Cells.Item(709,11).NumberFormat = "0.0000"
$ws.Cells.Item(709,12).Formula = '=0.49*0.4'
$ws.Cells.Item(709,13).Formula = '=2.5/G709'
$ws.Cells.Item(709,14).Value = 0
$ws.Cells.Item(709,15).Value = 0
$ws.Cells.Item(709,16).Value = 0
$ws.Cells.Item(709,17).Value = 0
$ws.Cells.Item(709,22).Value = 0

# Blank trailing rows with only K number-format set (style 6)
$ws.Cells.Item(710,11).NumberFormat = "0.0000"
$ws.Cells.Item(711,11).NumberFormat = "0.0000"
$ws.Cells.Item(712,11).NumberFormat = "0.0000"
$ws.Cells.Item(713,11).NumberFormat = "0.0000"
$ws.Cells.Item(714,11).NumberFormat = "0.0000"

# Match final selection / active cell seen in the authored workbook
$ws.Range("X707").Select()
